$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.725.18'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.39%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.290.55'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +3.46%  '

# Row 4
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.87'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.11%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.638'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.53%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.07'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +8.97%  '

# Row 8
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.645'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +3.73%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.37'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.38%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0990'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +5.43%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '59.22'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.37%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.30'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +3.53%  '

# Row 14
$ws.Range("E14").Value = '  +1.37%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.633.84'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +3.56%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.42'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +6.05%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.879'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.82%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.288.15'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.56%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.663.91'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.61%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0998'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.75%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.30'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.45%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.53'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.25%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.46'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.63%  '

# Row 24
$ws.Range("E24").Value = '  +10.17%  '

# Row 25
$ws.Range("E25").Value = '  +0.77%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.44'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.16%  '

# Row 27
$ws.Range("E27").Value = '  -0.14%  '

# Row 28
$ws.Range("E28").Value = '  +0.18%  '

# Row 29
$ws.Range("E29").Value = '  -1.77%  '

# Row 30
$ws.Range("E30").Value = '  -0.46%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '166.98'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.19%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.08'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.18%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.42'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +9.22%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.127'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +4.63%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0820'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +4.94%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '32.24'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +24.35%  '

# Row 37
$ws.Range("E37").Value = '  +3.62%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.73'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +15.57%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.76'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.54%  '

# Row 40
$ws.Range("E40").Value = '  -0.83%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.46'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +19.77%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.34'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +5.03%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.96'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +5.17%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.213'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +8.63%  '

# Row 45
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.13'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +6.83%  '

# Row 46
$ws.Range("B46").Value = 'MultiversX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.99'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.67%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.86'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -6.92%  '

# Row 48
$ws.Range("E48").Value = '  +3.52%  '

# Row 49
$ws.Range("E49").Value = '  +0.12%  '

# Row 50
$ws.Range("E50").Value = '  +2.60%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '98.36'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +5.73%  '
